$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

$ws.Range("E2").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F2").Value = "Ortega Valle Manuel"
$ws.Range("E4").Value = "INGLÉS IV"
$ws.Range("F4").Value = "González Nuñez Veronica"
$ws.Range("E5").Value = "FÍSICA I"
$ws.Range("F5").Value = "Polanco Domínguez Rosa María"
$ws.Range("E6").Value = "ECOLOGÍA"
$ws.Range("F6").Value = "Camarillo Aburto Raymundo"
$ws.Range("E7").Value = "FÍSICA I"
$ws.Range("F7").Value = "Polanco Domínguez Rosa María"
$ws.Range("E8").Value = "INGLÉS IV"
$ws.Range("F8").Value = "González Nuñez Veronica"
$ws.Range("E9").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F9").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E10").Value = "ECOLOGÍA"
$ws.Range("F10").Value = "Camarillo Aburto Raymundo"
$ws.Range("E13").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F13").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E14").Value = "INGLÉS IV"
$ws.Range("F14").Value = "González Nuñez Veronica"
$ws.Range("E15").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F15").Value = "Ortega Valle Manuel"
$ws.Range("E17").Value = "ECOLOGÍA"
$ws.Range("F17").Value = "Camarillo Aburto Raymundo"
$ws.Range("E18").Value = "INGLÉS IV"
$ws.Range("F18").Value = "González Nuñez Veronica"
$ws.Range("E20").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F20").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E21").Value = "FÍSICA I"
$ws.Range("F21").Value = "Polanco Domínguez Rosa María"
$ws.Range("E22").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F22").Value = "Ortega Valle Manuel"
$ws.Range("E23").Value = "FÍSICA I"
$ws.Range("F23").Value = "Polanco Domínguez Rosa María"
$ws.Range("E24").Value = "ECOLOGÍA"
$ws.Range("F24").Value = "Camarillo Aburto Raymundo"
$ws.Range("E26").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F26").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E27").Value = "INGLÉS IV"
$ws.Range("F27").Value = "González Nuñez Veronica"
$ws.Range("E29").Value = "ECOLOGÍA"
$ws.Range("F29").Value = "Camarillo Aburto Raymundo"
$ws.Range("E30").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F30").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E31").Value = "FÍSICA I"
$ws.Range("F31").Value = "Polanco Domínguez Rosa María"
$ws.Range("E34").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F34").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E35").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F35").Value = "Ortega Valle Manuel"
$ws.Range("E39").Value = "FÍSICA I"
$ws.Range("F39").Value = "Polanco Domínguez Rosa María"
$ws.Range("E40").Value = "ECOLOGÍA"
$ws.Range("F40").Value = "Camarillo Aburto Raymundo"
$ws.Range("E42").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F42").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E43").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F43").Value = "Ortega Valle Manuel"
$ws.Range("E45").Value = "INGLÉS IV"
$ws.Range("F45").Value = "González Nuñez Veronica"
$ws.Range("E46").Value = "FÍSICA I"
$ws.Range("F46").Value = "Polanco Domínguez Rosa María"
$ws.Range("E47").Value = "ECOLOGÍA"
$ws.Range("F47").Value = "Camarillo Aburto Raymundo"
$ws.Range("E48").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F48").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E49").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F49").Value = "Ortega Valle Manuel"
$ws.Range("E50").Value = "ECOLOGÍA"
$ws.Range("F50").Value = "Camarillo Aburto Raymundo"
$ws.Range("E52").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F52").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E54").Value = "FÍSICA I"
$ws.Range("F54").Value = "Polanco Domínguez Rosa María"
$ws.Range("E55").Value = "ECOLOGÍA"
$ws.Range("F55").Value = "Camarillo Aburto Raymundo"
$ws.Range("E56").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F56").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E57").Value = "INGLÉS IV"
$ws.Range("F57").Value = "González Nuñez Veronica"
$ws.Range("E58").Value = "FÍSICA I"
$ws.Range("F58").Value = "Polanco Domínguez Rosa María"
$ws.Range("E59").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F59").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E60").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F60").Value = "Ortega Valle Manuel"
$ws.Range("E61").Value = "ECOLOGÍA"
$ws.Range("F61").Value = "Camarillo Aburto Raymundo"
$ws.Range("E62").Value = "INGLÉS IV"
$ws.Range("F62").Value = "González Nuñez Veronica"
$ws.Range("E63").Value = "INGLÉS IV"
$ws.Range("F63").Value = "González Nuñez Veronica"
$ws.Range("E64").Value = "FÍSICA I"
$ws.Range("F64").Value = "Polanco Domínguez Rosa María"
$ws.Range("E65").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F65").Value = "Ortega Valle Manuel"
$ws.Range("E66").Value = "ECOLOGÍA"
$ws.Range("F66").Value = "Camarillo Aburto Raymundo"
$ws.Range("E70").Value = "FÍSICA I"
$ws.Range("F70").Value = "Polanco Domínguez Rosa María"
$ws.Range("E72").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F72").Value = "Ortega Valle Manuel"
$ws.Range("E73").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F73").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E74").Value = "ECOLOGÍA"
$ws.Range("F74").Value = "Camarillo Aburto Raymundo"
$ws.Range("E80").Value = "INGLÉS IV"
$ws.Range("F80").Value = "González Nuñez Veronica"
$ws.Range("E81").Value = "ECOLOGÍA"
$ws.Range("F81").Value = "Camarillo Aburto Raymundo"
$ws.Range("E82").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F82").Value = "Ortega Valle Manuel"
$ws.Range("E83").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F83").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E85").Value = "ECOLOGÍA"
$ws.Range("F85").Value = "Camarillo Aburto Raymundo"
$ws.Range("E90").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F90").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E91").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F91").Value = "Ortega Valle Manuel"
$ws.Range("E92").Value = "FÍSICA I"
$ws.Range("F92").Value = "Polanco Domínguez Rosa María"
$ws.Range("E93").Value = "INGLÉS IV"
$ws.Range("F93").Value = "González Nuñez Veronica"
$ws.Range("E95").Value = "INGLÉS IV"
$ws.Range("F95").Value = "González Nuñez Veronica"
$ws.Range("E96").Value = "FÍSICA I"
$ws.Range("F96").Value = "Polanco Domínguez Rosa María"
$ws.Range("E97").Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws.Range("F97").Value = "Zarate Amezcua Eladio Jorge"
$ws.Range("E98").Value = "ECOLOGÍA"
$ws.Range("F98").Value = "Camarillo Aburto Raymundo"
$ws.Range("E100").Value = "CÁLCULO DIFERENCIAL"
$ws.Range("F100").Value = "Ortega Valle Manuel"
$ws.Range("E101").Value = "FÍSICA I"
$ws.Range("F101").Value = "Polanco Domínguez Rosa María"
$ws.Range("E103").Value = "ECOLOGÍA"
$ws.Range("F103").Value = "Camarillo Aburto Raymundo"
